$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before the current row 3 (4x5090) to make room for the new entries
$ws.Rows.Item(3).Resize(3).EntireRow.Insert()

# Row 3: 2x4090
$ws.Cells.Item(3, 1).Value = "2x4090"
$ws.Cells.Item(3, 2).Value = 467.26
$ws.Cells.Item(3, 3).Value = 0.78
$ws.Cells.Item(3, 4).Value = 0.4636961577422991

# Row 4: 2x5090
$ws.Cells.Item(4, 1).Value = "2x5090"
$ws.Cells.Item(4, 2).Value = 1230.14
$ws.Cells.Item(4, 3).Value = 1.3
$ws.Cells.Item(4, 4).Value = 0.2935528566757533

# Row 5: 4x4090
$ws.Cells.Item(5, 1).Value = "4x4090"
$ws.Cells.Item(5, 2).Value = 906.1900000000001
$ws.Cells.Item(5, 3).Value = 1.56
$ws.Cells.Item(5, 4).Value = 0.4781925791868519

# Row 6: 4x5090 (updated values)
$ws.Cells.Item(6, 2).Value = 2501.38
$ws.Cells.Item(6, 4).Value = 0.2887295101992589
